# Apply "parts list" update:
#  - C4: big bypass capacitor changed from Tantalum to Ceramic part
#  - C13-C17: changed from 47nF to 6800pF (also X5R -> X7R) ceramic part
#  - R7-R11: changed from 1.5MOhm to 10MOhm resistor part
#  - U2 description simplified ("2.5 V low dropout linear voltage regulator" -> "2.5 V linear voltage regulator")
#  - G3 (U1 Digikey link) gets turned into a real hyperlink
#  - Selection moves to A8
#  - Columns A, D, E, F widths adjust to fit the new (longer) content

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- R7-R11 (row 16): 0201 1.5MOhm -> 0201 10MOhm (same manufacturer, Rohm) ---
$ws.Range("D16").Value = "0201, 10MOhm, 1%, 1/20W"
$ws.Range("F16").Value = "MCR006YRTF1005"
$ws.Range("G16").Value = "http://www.digikey.com/short/38853w"

# --- C13-C17 (row 13): 0201 Ceramic X5R 47nF (Murata) -> 0201 Ceramic X7R 6800pF (TDK) ---
$ws.Range("E13").Value = "TDK Corporation"
$ws.Range("F13").Value = "CGA1A2X7R1C682K030BA"
$ws.Range("G13").Value = "http://www.digikey.com/short/38852m"
$ws.Range("D13").Value = "0201, Ceramic, X7R, 6800pF, 10%, 16V"

# --- U2 (row 4): simplify description text ---
$ws.Range("D4").Value = "2.5 V linear voltage regulator"

# --- C4 (row 8): 10uF Tantalum (AVX) -> 0502 Ceramic X5R 10uF (Samsung) ---
$ws.Range("E8").Value = "Samsung Electro-Mechanics America, Inc."
$ws.Range("D8").Value = "0502, Ceramic, X5R, 10uF, 20% 10V"
$ws.Range("F8").Value = "CL05A106MP5NUNC"
$ws.Range("G8").Value = "http://www.digikey.com/short/388517"

# --- G3: make the U1 Digikey link a real hyperlink (it already shows the URL text) ---
$ws.Hyperlinks.Add($ws.Range("G3"), "http://www.digikey.com/short/3ttfj7")
$ws.Range("G3").Style = "Hyperlink"

# --- column widths: re-fit the columns whose longest entry changed ---
$ws.Columns.Item(1).ColumnWidth = 13.833333333333334
$ws.Columns.Item(4).ColumnWidth = 33.666666666666664
$ws.Columns.Item(5).ColumnWidth = 37.5
$ws.Columns.Item(6).ColumnWidth = 22.666666666666668

# --- selection moves to A8 ---
$ws.Range("A8").Select()
